$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab to reflect the new "through" date
$ws.Name = "Through 2022-06-22"

# Update the column header text (shared string) for the current-month column
$ws.Range("B1").Value = "June 2022 (through June 22)"

# Update existing cell values that changed
$ws.Range("H3").Value = 2
$ws.Range("H4").Value = 5
$ws.Range("N4").Value = 3
$ws.Range("AF6").Value = 2
$ws.Range("Z10").Value = 3
$ws.Range("AF10").Value = 3
$ws.Range("AL10").Value = 4
$ws.Range("H14").Value = 10
$ws.Range("AL14").Value = 1
$ws.Range("AR14").Value = 2
$ws.Range("B15").Value = 3
$ws.Range("H17").Value = 2
$ws.Range("Z26").Value = 1
$ws.Range("B38").Value = 1
$ws.Range("N38").Value = 4
$ws.Range("H39").Value = 2
$ws.Range("AF41").Value = 3
$ws.Range("Z56").Value = 2
$ws.Range("B71").Value = 2
$ws.Range("H92").Value = 2
